$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "tomastheBest"

$ws.Range("A7").Select()
